$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''64.664.76'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.94%  '

$ws.Range('D3').Value = '''3.137.72'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.03%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '''571.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.20%  '

$ws.Range('D6').Value = '''147.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.22%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').Value = '''3.136.30'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.14%  '

$ws.Range('D9').Value = '''0.522'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.73%  '

$ws.Range('E10').Value = '  -3.99%  '

$ws.Range('D11').Value = '''6.05'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.30%  '

$ws.Range('D12').Value = '''0.495'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.78%  '

$ws.Range('E13').Value = '  +0.42%  '

$ws.Range('D14').Value = '''36.74'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.11%  '

$ws.Range('D15').Value = '''3.652.59'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.19%  '

$ws.Range('D16').Value = '''64.830.51'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.55%  '

$ws.Range('D17').Value = '''3.136.72'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.50%  '

$ws.Range('D18').Value = '''7.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.54%  '

$ws.Range('D19').Value = '''0.111'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.21%  '

$ws.Range('D20').Value = '''498.54'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.26%  '

$ws.Range('D21').Value = '''14.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.52%  '

$ws.Range('D22').Value = '''0.708'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.19%  '

$ws.Range('E23').Value = '  -3.69%  '

$ws.Range('D24').Value = '''7.64'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.03%  '

$ws.Range('D25').Value = '''83.48'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.15%  '

$ws.Range('E26').Value = '  -0.43%  '

$ws.Range('D27').Value = '''2.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.64%  '

$ws.Range('D28').Value = '''8.78'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.23%  '

$ws.Range('D29').Value = '''2.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.84%  '

$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '''27.29'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.39%  '

$ws.Range('B31').Value = 'Stacks'
$ws.Range('C31').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D31').Value = '''2.75'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.97%  '

$ws.Range('E32').Value = '  -0.13%  '

$ws.Range('E33').Value = '  -0.70%  '

$ws.Range('D34').Value = '''6.12'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.18%  '

$ws.Range('D35').Value = '''6.40'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.80%  '

$ws.Range('D36').Value = '''54.41'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.36%  '

$ws.Range('D37').Value = '''0.0893'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.63%  '

$ws.Range('D38').Value = '''467.68'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.10%  '

$ws.Range('D39').Value = '''0.0411'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.84%  '

$ws.Range('E40').Value = '  -1.99%  '

$ws.Range('D41').Value = '''8.58'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.17%  '

$ws.Range('D42').Value = '''3.011.58'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.69%  '

$ws.Range('E43').Value = '  -4.88%  '

$ws.Range('D44').Value = '''0.280'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.98%  '

$ws.Range('E45').Value = '  -2.33%  '

$ws.Range('D46').Value = '''28.02'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.03%  '

$ws.Range('D47').Value = '''0.0₃0570'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.18%  '

$ws.Range('E49').Value = '  -2.00%  '

$ws.Range('D50').Value = '''2.22'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.86%  '

$ws.Range('D51').Value = '''117.43'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.63%  '
